$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like numbers (contain only digits/dot) need to be
# forced to Text format first, otherwise Excel auto-converts them to numeric values
# and rounds/reformats them (e.g. "247.82" -> 247.81999999999999).
# NumberFormat is reset back to the default ("Normal" style) afterwards so the
# saved cell carries no explicit style, matching the original file layout.
$numericLookingCells = @(
    'D5', 'D7', 'D8', 'D10', 'D11', 'D12', 'D14', 'D15', 'D16', 'D20', 'D22', 'D23', 'D24', 'D25', 'D26', 'D28', 'D29', 'D30', 'D32', 'D34', 'D36', 'D38', 'D39', 'D40', 'D41', 'D42', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51'
)
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply all the updated values (coin prices + volume percentages), and the
# Aptos/Elrond row swap (rows 48-49).
$ws.Range('D2').Value = '30.551.19'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '1.873.76'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '247.82'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D7').Value = '0.4735'
$ws.Range('E7').Value = '  -0.98%  '
$ws.Range('D8').Value = '0.2898'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -1.31%  '
$ws.Range('D10').Value = '21.99'
$ws.Range('E10').Value = '  +2.69%  '
$ws.Range('D11').Value = '0.07738'
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').Value = '0.7422'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = '1.875.15'
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('D14').Value = '96.03'
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').Value = '5.168'
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('D16').Value = '274.91'
$ws.Range('E16').Value = '  -1.85%  '
$ws.Range('D17').Value = '30.625.82'
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('E18').Value = '  -3.04%  '
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').Value = '0.000007469'
$ws.Range('E20').Value = '  -2.13%  '
$ws.Range('D21').Value = '2.118.84'
$ws.Range('E21').Value = '  -0.85%  '
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').Value = '5.203'
$ws.Range('E23').Value = '  -2.10%  '
$ws.Range('D24').Value = '6.167'
$ws.Range('E24').Value = '  -1.29%  '
$ws.Range('D25').Value = '9.192'
$ws.Range('E25').Value = '  -1.56%  '
$ws.Range('D26').Value = '164.94'
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('E27').Value = '  -2.36%  '
$ws.Range('D28').Value = '1.903'
$ws.Range('E28').Value = '  -4.92%  '
$ws.Range('D29').Value = '0.09947'
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('D30').Value = '1.347'
$ws.Range('E30').Value = '  -2.32%  '
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('D32').Value = '4.234'
$ws.Range('E32').Value = '  -2.88%  '
$ws.Range('E33').Value = '  -1.45%  '
$ws.Range('D34').Value = '0.04770'
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('E35').Value = '  -1.47%  '
$ws.Range('D36').Value = '0.6922'
$ws.Range('E36').Value = '  -1.99%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = '0.01848'
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('D39').Value = '2.753'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').Value = '6.265'
$ws.Range('E40').Value = '  -4.29%  '
$ws.Range('D41').Value = '73.21'
$ws.Range('E41').Value = '  +3.04%  '
$ws.Range('D42').Value = '1.970'
$ws.Range('E42').Value = '  +2.12%  '
$ws.Range('D44').Value = '0.4159'
$ws.Range('E44').Value = '  -1.10%  '
$ws.Range('D45').Value = '0.8336'
$ws.Range('E45').Value = '  -2.14%  '
$ws.Range('D46').Value = '101.10'
$ws.Range('E46').Value = '  -1.71%  '
$ws.Range('D47').Value = '9.345'
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').Value = '35.33'
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').Value = '6.973'
$ws.Range('E49').Value = '  -2.73%  '
$ws.Range('D50').Value = '911.93'
$ws.Range('E50').Value = '  -2.09%  '
$ws.Range('D51').Value = '0.05665'

# Restore default (un-styled) formatting on the cells we forced to Text above.
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).Style = "Normal"
}
